$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The watch-list table in B:F (rows 2-19) is replaced with a refreshed
# symbol list, and the table grows from 19 to 32 rows (A column keeps
# counting 0,1,2,...).
#
# Each inner array is: @(rowNumber, B, C, D, E, F) - use $null for a
# cell that should end up blank.
$data = @(
    @(2, "NSE:AFFLE", "NSE:3IINFOLTD", "NSE:NAUKRI", "NSE:ABBOTINDIA", "NSE:ASTRAL"),
    @(3, "NSE:AJMERA", "NSE:EIDPARRY", $null, "NSE:ACC", "NSE:AUBANK"),
    @(4, "NSE:ALPA", "NSE:GIPCL", $null, "NSE:BPCL", "NSE:CUMMINSIND"),
    @(5, "NSE:ASAHIINDIA", "NSE:GRAPHITE", $null, "NSE:COALINDIA", "NSE:NAUKRI"),
    @(6, "NSE:ASTRAL", "NSE:HEALTHY", $null, "NSE:IRCTC", "NSE:PAGEIND"),
    @(7, "NSE:AVALON", "NSE:IFCI", $null, $null, $null),
    @(8, "NSE:BAJAJHCARE", "NSE:IRISDOREME", $null, $null, $null),
    @(9, "NSE:BBTC", "NSE:KARURVYSYA", $null, $null, $null),
    @(10, "NSE:BOMDYEING", "NSE:MAGNUM", $null, $null, $null),
    @(11, "NSE:CGPOWER", "NSE:MANAKSIA", $null, $null, $null),
    @(12, "NSE:CUMMINSIND", "NSE:RAMAPHO", $null, $null, $null),
    @(13, "NSE:DCXINDIA", "NSE:RKEC", $null, $null, $null),
    @(14, "NSE:DVL", $null, $null, $null, $null),
    @(15, "NSE:ENIL", $null, $null, $null, $null),
    @(16, "NSE:GANECOS", $null, $null, $null, $null),
    @(17, "NSE:GOCLCORP", $null, $null, $null, $null),
    @(18, "NSE:GPPL", $null, $null, $null, $null),
    @(19, "NSE:GULFOILLUB", $null, $null, $null, $null),
    @(20, "NSE:HINDWAREAP", $null, $null, $null, $null),
    @(21, "NSE:JINDALPOLY", $null, $null, $null, $null),
    @(22, "NSE:KICL", $null, $null, $null, $null),
    @(23, "NSE:NAUKRI", $null, $null, $null, $null),
    @(24, "NSE:NEULANDLAB", $null, $null, $null, $null),
    @(25, "NSE:NUCLEUS", $null, $null, $null, $null),
    @(26, "NSE:OMAXE", $null, $null, $null, $null),
    @(27, "NSE:PAGEIND", $null, $null, $null, $null),
    @(28, "NSE:PCBL", $null, $null, $null, $null),
    @(29, "NSE:PDMJEPAPER", $null, $null, $null, $null),
    @(30, "NSE:PIXTRANS", $null, $null, $null, $null),
    @(31, "NSE:PYRAMID", $null, $null, $null, $null),
    @(32, "NSE:ROTO", $null, $null, $null, $null)
)

# Rows 20-32 are brand new. Copy the formatting (border/alignment/font)
# from the existing styled column-A cell A19 down onto A20:A32 first, so
# the new index cells pick up the same style used by the rest of column A.
$ws.Range("A19").Copy()
$ws.Range("A20:A32").PasteSpecial(-4122)
$excel.CutCopyMode = 0

foreach ($r in $data) {
    $row = $r[0]
    $ws.Cells.Item($row, 1).Value = $row - 2
    for ($col = 2; $col -le 6; $col++) {
        $val = $r[$col - 1]
        if ($null -eq $val) { $val = "" }
        $ws.Cells.Item($row, $col).Value = $val
    }
}
